# Add kick off for second catalogues test
# -----------------------------------------------------------------
# This script reproduces, via the Excel COM object model, the edits
# made to AdminInputData.xlsx: new header columns and a second
# "AddNewCatalogue" test row on the CatalogueManagement sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Admin")
$ws2 = $wb.Worksheets.Item("CatalogueManagement")

# ---------------------------------------------------------------
# 1. New header cells on row 1 (F1:J1 new, E1 gets a value)
# ---------------------------------------------------------------
$ws2.Range("E1").Value = "shortName"

$ws2.Range("F1").Value = "longName"
$ws2.Range("F1").NumberFormat = "@"

$ws2.Range("G1").Value = "introDateDaysJump"
$ws2.Range("G1").NumberFormat = "@"

$ws2.Range("H1").Value = "exitDateDaysJump"
$ws2.Range("H1").NumberFormat = "@"

$ws2.Range("I1").Value = "configurationSet"
$ws2.Range("I1").NumberFormat = "@"

$ws2.Range("J1").Value = "description"
$ws2.Range("J1").NumberFormat = "@"

# ---------------------------------------------------------------
# 2. Row 3 : fill in the new "n/a" cells (E3:J3)
# ---------------------------------------------------------------
$ws2.Range("E3").Value = "n/a"

$ws2.Range("F3").Value = "n/a"
$ws2.Range("F3").NumberFormat = "@"
$ws2.Range("F3").Font.Bold = $true

$ws2.Range("G3").Value = "n/a"
$ws2.Range("G3").NumberFormat = "@"
$ws2.Range("G3").Font.Bold = $true

$ws2.Range("H3").Value = "n/a"
$ws2.Range("H3").NumberFormat = "@"
$ws2.Range("H3").Font.Bold = $true

$ws2.Range("I3").Value = "n/a"
$ws2.Range("I3").NumberFormat = "@"
$ws2.Range("I3").Font.Bold = $true

$ws2.Range("J3").Value = "n/a"
$ws2.Range("J3").NumberFormat = "@"
$ws2.Range("J3").Font.Bold = $true

# ---------------------------------------------------------------
# 3. New row 5 : second catalogue test (AddNewCatalogue)
# ---------------------------------------------------------------
$ws2.Range("A5").Value = "AddNewCatalogue"
$ws2.Range("A5").NumberFormat = "GENERAL"
$ws2.Range("A5").Font.Name = "Century Gothic"
$ws2.Range("A5").Font.Size = 10.5
$ws2.Range("A5").Font.Color = 0

$ws2.Range("B5").Value = "positive"
$ws2.Range("C5").Value = "Nazar_Lelyak"

$ws2.Range("D5").Value = "qwerty123"
$ws2.Range("D5").NumberFormat = "GENERAL"
$ws2.Range("D5").HorizontalAlignment = -4131

$ws2.Range("E5").Value = "auto_test_catalogue"
$ws2.Range("F5").Value = "auto_test_catalogue_please_ignore"

$ws2.Rows.Item(5).RowHeight = 13.2

# ---------------------------------------------------------------
# 4. Column widths for the new columns (E, F, G)
# ---------------------------------------------------------------
$ws2.Columns.Item(5).ColumnWidth = 20.285714285714285
$ws2.Columns.Item(6).ColumnWidth = 10.857142857142858
$ws2.Columns.Item(7).ColumnWidth = 14.571428571428571

# ---------------------------------------------------------------
# 5. Selection / view state
# ---------------------------------------------------------------
$ws2.Activate()
$ws2.Range("I3:J3").Select()

$ws1.Activate()
$u = $excel.Union($ws1.Range("I3:J3"), $ws1.Range("A3"))
$u.Select()
$ws1.Range("A3").Select()
